# Added events insertion functionality
# Shift every existing event's Date (column B) forward by one day to make
# room for newly inserted events.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value = $current + 1
    }
}
